$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.000012051522389810998
$ws.Range("A3").Value = 0.000009512608812656254
$ws.Range("F3").Value = 8.0
$ws.Range("A4").Value = 0.000001978260797841358
$ws.Range("F4").Value = 1.0
$ws.Range("A5").Value = 0.000000560652154035779
$ws.Range("E5").Value = 10.0
$ws.Range("F5").Value = 6.0
